$p = $ppt.ActivePresentation

# Slide 1: Title "Header" + " " + "with" + " " + "inline code"(Courier)
#   -> "Header " + "with " + "inline code"(Courier)
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(1, 7).Text = "Header "
$tr1.Characters(8, 5).Text = "with "

# Slide 2: Title "Syntax" + " " + "highlighting"
#   -> "Syntax " + "highlighting"
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Characters(1, 7).Text = "Syntax "

# Slide 3: Title "Two" + " " + "column" + " " + "slide"
#   -> "Two " + "column " + "slide"
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Characters(1, 4).Text = "Two "
$tr3.Characters(5, 7).Text = "column "
